$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the existing two data rows (currently sitting in C1:N2) ---
$oldRows = @(
    @($ws.Cells.Item(1,3).Value2, $ws.Cells.Item(1,4).Value2, $ws.Cells.Item(1,5).Value2, $ws.Cells.Item(1,6).Value2, $ws.Cells.Item(1,7).Value2, $ws.Cells.Item(1,8).Value2, $ws.Cells.Item(1,9).Value2, $ws.Cells.Item(1,10).Value2, $ws.Cells.Item(1,11).Value2, $ws.Cells.Item(1,12).Value2, $ws.Cells.Item(1,13).Value2, $ws.Cells.Item(1,14).Value2),
    @($ws.Cells.Item(2,3).Value2, $ws.Cells.Item(2,4).Value2, $ws.Cells.Item(2,5).Value2, $ws.Cells.Item(2,6).Value2, $ws.Cells.Item(2,7).Value2, $ws.Cells.Item(2,8).Value2, $ws.Cells.Item(2,9).Value2, $ws.Cells.Item(2,10).Value2, $ws.Cells.Item(2,11).Value2, $ws.Cells.Item(2,12).Value2, $ws.Cells.Item(2,13).Value2, $ws.Cells.Item(2,14).Value2)
)

# Wipe the old range entirely (values + formats) so nothing stale survives the move.
$ws.Range("C1:N2").Clear()

# --- Re-write the two data rows, shifted left two columns (C->A) and down one row (1->2) ---
for ($r = 0; $r -lt 2; $r++) {
    $destRow = $r + 2
    for ($c = 0; $c -lt 12; $c++) {
        $destCol = $c + 1
        $ws.Cells.Item($destRow, $destCol).Value = $oldRows[$r][$c]
    }
}

# Re-apply the numeric "0" format to the numeric columns (D,E,F,G,I,J,K,L) on both rows,
# matching the style used before the move.
$numCols = @(4,5,6,7,9,10,11,12)
foreach ($col in $numCols) {
    $ws.Cells.Item(2, $col).NumberFormat = "0"
    $ws.Cells.Item(3, $col).NumberFormat = "0"
}

# --- Insert the new header row on row 1 ---
# The 10 columns that already existed in the source data get their header
# typed first (left to right, skipping the two brand-new columns), then the
# two newly-inserted columns (JENIS_PRODUK, SUKU_BUNGA_EFFECTIVE) are filled
# in last.
$ws.Cells.Item(1, 1).Value = "PARTNER"
$ws.Cells.Item(1, 2).Value = "DEBITUR"
$ws.Cells.Item(1, 4).Value = "NILAI_PEMBIAYAAN_POKOK_MAXIMUM"
$ws.Cells.Item(1, 5).Value = "SUKU_BUNGA_FLAT"
$ws.Cells.Item(1, 7).Value = "JANGKA_WAKTU_MAXIMUM"
$ws.Cells.Item(1, 8).Value = "POLA_PEMBAYARAN"
$ws.Cells.Item(1, 9).Value = "BIAYA_ADMINISTRASI"
$ws.Cells.Item(1, 10).Value = "BIAYA_ASURANSI"
$ws.Cells.Item(1, 11).Value = "BIAYA_PROVINSI"
$ws.Cells.Item(1, 12).Value = "BIAYA_LAIN_LAIN"
$ws.Cells.Item(1, 3).Value = "JENIS_PRODUK"
$ws.Cells.Item(1, 6).Value = "SUKU_BUNGA_EFFECTIVE"
$ws.Range("A1:L1").Font.Bold = $true

# --- Turn the range into an Excel Table ---
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:L3"), $null, 1)
$lo.Name = "Table2"
$lo.TableStyle = ""
$lo.ShowTableStyleRowStripes = $false
$lo.ShowTableStyleColumnStripes = $false
$lo.ShowTableStyleFirstColumn = $false
$lo.ShowTableStyleLastColumn = $false
$lo.ShowAutoFilterDropDown = $false

# --- Sheet view: freeze header row, restore the recorded selections ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("G1").Select()
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("F2").Select()
